$wb = $excel.ActiveWorkbook

# --- Sheet "Info" (sheet1): update Objetivo / Tiempo summary values ---
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Cells.Item(2, 1).Value = 640108574274.0112
$wsInfo.Cells.Item(2, 2).Value = 2.092000007629395

# --- Sheet "Activados" (sheet2): Proceso changes to 1, Tiempo extended to 20 rows (step 20) ---
$wsActivados = $wb.Worksheets.Item("Activados")
for ($r = 2; $r -le 20; $r++) {
    $wsActivados.Cells.Item($r, 1).Value = 1
    $wsActivados.Cells.Item($r, 2).Value = ($r - 2) * 20
}

# --- Sheet "Operando" (sheet3): Proceso column changes from 4 to 1 for rows 2..366 ---
$wsOperando = $wb.Worksheets.Item("Operando")
for ($r = 2; $r -le 366; $r++) {
    $wsOperando.Cells.Item($r, 1).Value = 1
}

# --- Sheet "Contaminantes" (sheet6): update mass / concentration values ---
$wsContaminantes = $wb.Worksheets.Item("Contaminantes")
$wsContaminantes.Cells.Item(2, 2).Value = 449208244800.0004
$wsContaminantes.Cells.Item(2, 3).Value = 16.66000000000001
$wsContaminantes.Cells.Item(3, 2).Value = 13481640000.00001
$wsContaminantes.Cells.Item(3, 3).Value = 0.5000000000000004
$wsContaminantes.Cells.Item(4, 2).Value = 87091394399.99998
$wsContaminantes.Cells.Item(4, 3).Value = 3.23
$wsContaminantes.Cells.Item(5, 2).Value = 307074.010608
$wsContaminantes.Cells.Item(5, 3).Value = 0.0000113886
$wsContaminantes.Cells.Item(6, 2).Value = 90326988000.00008
$wsContaminantes.Cells.Item(6, 3).Value = 3.350000000000003
